$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -2
    "F4"  = -1
    "F5"  = 4
    "F7"  = -3
    "F20" = 2
    "F30" = -1
    "F39" = 1
    "F40" = 1
    "F50" = 2
    "F51" = 3
    "F59" = 0
    "F75" = 0
    "F76" = 0
    "F84" = 0
    "F86" = 5
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
